$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: Mock up - Native APP, 1h30 ---
$ws.Range("A3").Value = 41254
$ws.Range("A3").NumberFormat = "mm-dd-yy"
$ws.Range("B3").Value = "Mock up - Native APP"
$ws.Range("D3").Value = "1h30"

# --- Row 4: Logo / LogoComponents, 3h ---
$ws.Range("A4").Value = 41262
$ws.Range("B4").Value = "Logo"
$ws.Range("C4").Value = "LogoComponents"
$ws.Range("D4").Value = "3h"

# --- Row 5: Logo / Logo zelf, 1h ---
$ws.Range("A5").Value = 41262
$ws.Range("B5").Value = "Logo"
$ws.Range("C5").Value = "Logo zelf"
$ws.Range("D5").Value = "1h"

# --- Row 6: Logo / Logo zelf, 3h30 ---
$ws.Range("A6").Value = 41263
$ws.Range("B6").Value = "Logo"
$ws.Range("C6").Value = "Logo zelf"
$ws.Range("D6").Value = "3h30"

# --- Row 7: Moodboard ---
$ws.Range("A7").Value = 41263
$ws.Range("B7").Value = "Moodboard"

# Give A4:A7 the same date number format/style as A3 (reuse one style record)
$ws.Range("A3").Copy()
$ws.Range("A4:A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Selection moves to D7 (single cell)
[void]$ws.Range("D7").Select()
